# "library , IQAC , Staff excel are updated"
# Staff name column (A) gets honorific prefixes (Mr./Mrs./Dr) added, and the
# Balamurugapandian row is fully upper-cased. Column A is also widened to fit
# the longer names, and the saved selection/scroll state is reset to A12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update staff names with honorifics (order matches the edit session that
# produced the shared-string table in the target workbook).
$ws.Range("A9").Value  = "Dr BALAMURUGAPANDIAN N"
$ws.Range("A3").Value  = "Mrs. SASIKALA S"
$ws.Range("A4").Value  = "Mrs. RAMA NACHIAR R"
$ws.Range("A6").Value  = "Mr. SIVARAJ C"
$ws.Range("A7").Value  = "Mrs. KAVITHA K"
$ws.Range("A8").Value  = "Mr. NIRUBAN BALU T"
$ws.Range("A11").Value = "Mrs. MOUNICA A"
$ws.Range("A12").Value = "Mr. KUMAR V"

# Widen column A (Name) now that values are longer, and drop the old
# "best fit" auto-width in favor of an explicit custom width.
$ws.Columns.Item(1).ColumnWidth = 27

# Reset view state: clear the scrolled-away top-left cell and move the
# selection to A12.
[void]$ws.Range("A12").Select()
